$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Footer (Primary) -- contains the Pearson logo whose wp:docPr/@id="2" ---
# Rename it from "image1.png" to "image2.png" (wp:docPr name + pic:cNvPr name).
$footerPrimary = $sec.Footers(1)
if ($footerPrimary.Exists) {
    $shapeCount = $footerPrimary.Range.InlineShapes.Count
    for ($i = 1; $i -le $shapeCount; $i++) {
        $ish = $footerPrimary.Range.InlineShapes($i)
        if ($ish.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shp = $ish.ConvertToShape()
            $shp.Name = "image2.png"
            $shp.ConvertToInlineShape()
        }
    }
}

# --- Footer (First Page) -- contains the Pearson logo whose wp:docPr/@id="3" ---
# Rename it from "image1.png" to "image2.png" (wp:docPr name + pic:cNvPr name).
$footerFirstPage = $sec.Footers(2)
if ($footerFirstPage.Exists) {
    $shapeCount = $footerFirstPage.Range.InlineShapes.Count
    for ($i = 1; $i -le $shapeCount; $i++) {
        $ish = $footerFirstPage.Range.InlineShapes($i)
        if ($ish.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shp = $ish.ConvertToShape()
            $shp.Name = "image2.png"
            $shp.ConvertToInlineShape()
        }
    }
}

# --- Header (First Page) -- contains the BTEC logo (wp:docPr/@id="1") ---
# Rename it from "image2.jpg" to "image1.jpg" (wp:docPr name + pic:cNvPr name).
$headerFirstPage = $sec.Headers(2)
if ($headerFirstPage.Exists) {
    $shapeCount = $headerFirstPage.Range.InlineShapes.Count
    for ($i = 1; $i -le $shapeCount; $i++) {
        $ish = $headerFirstPage.Range.InlineShapes($i)
        if ($ish.AlternativeText -eq "BTec_Logo-Orange") {
            $shp = $ish.ConvertToShape()
            $shp.Name = "image1.jpg"
            $shp.ConvertToInlineShape()
        }
    }
}
